$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: replacement failed for: $findText"
    }
    return $ok
}

# 1. Rewrite the paragraph describing the UK/Thai data scientists & astronomers'
#    respective track records (JE-S section).
Replace-Text `
    "The data scientists involved in the project have a strong track record in developing ML algorithms, including automated image analysis, (Boongoen, Eungwanichayapant, Iam-On, Uttuma) and setting-up and maintaining DM systems (Boongoen, Iam-On). The UK and Thai astronomers have track records in using non-ML techniques to analyse large astronomical datasets and are all involved in the Gravitational-wave Optical Transient Observatory (GOTO) project (see section 3 for a description of GOTO), which is the source of the large astronomical datasets that is being used throughout the project. Importantly, through our collaborative work during Phase 1 project, all data scientists and astronomers involved in the project now also have a good working knowledge of each other's areas of expertise (i.e., the astronomers have learned from the data scientists, and vice versa). This cross-disciplinary knowledge makes our team especially well-suited to exploiting astronomical datasets to train others in advanced data handling techniques, particularly those which a background in astronomy or other physical sciences." `
    "The Thai data science partners have a strong track record in developing ML algorithms, including automated image analysis, (Boongoen, Eungwanichayapant, Iam-On, Uttuma) and setting-up and maintaining DM systems (Boongoen, Iam-On). The UK partners have extensive understanding of the methods and outputs of the pipelines used to process data from the Gravitational-wave Optical Transient Observatory (GOTO; see section 3 for a description of GOTO), which is the source of the large astronomical datasets that are used throughout the project. The Thai astronomy partners (Sawangwit, Awiphan) are familiar with and have ready access to the computing infrastructure based at NARIT that will be used throughout the project. Importantly, through our collaborative work during Phase 1, all data scientists and astronomers involved in the project now also have a good working knowledge of each other's areas of expertise. This cross-disciplinary knowledge makes our team especially well-suited to exploiting astronomical datasets to train others – especially astronomers – advanced data handling techniques."

# 2. Clarify that AWS stands for Amazon Web Services.
Replace-Text `
    "leasing facilities provided by Amazon Web Services. This " `
    "leasing facilities provided by Amazon Web Services (AWS). This "

# 3. Secondary objective now refers to establishing a Thai data centre rather
#    than developing data handling systems.
Replace-Text `
    "A secondary objective is the development of data handling systems to organise and analyse" `
    "A secondary objective is the establishment of a Thai data centre to store, organise and analyse"

# 4. Student will work alongside "the UK partners" instead of "GOTO scientists".
Replace-Text `
    "will work alongside GOTO scientists with further guidance from Boongoen and Iam-On to set-up" `
    "will work alongside the UK partners with further guidance from Boongoen and Iam-On to set-up"

# 5. "major" -> "significant" research challenge.
Replace-Text `
    "enough to handle GOTO's data rate is a major research challenge." `
    "enough to handle GOTO's data rate is a significant research challenge."

# 6. Remove redundant "the" before "research" in the AWS/hardware sentence.
Replace-Text `
    "systems require different hardware infrastructure, the first stages of the research will be conducted on AWS" `
    "systems require different hardware infrastructure, the first stages of research will be conducted on AWS"

# 7. Clarify infrastructure will not be purchased "by NARIT" until optimum
#    system identified.
Replace-Text `
    "infrastructure until after the optimum system has been identified" `
    "infrastructure by NARIT until after the optimum system has been identified"

# 8. Note that "establishing" the data centre is itself an important research
#    project.
Replace-Text `
    "This research will make the Thai data centre itself an important" `
    "This research will make establishing the Thai data centre itself an important"

# 9. Update the cached PAGE field result in the footer from "5" to "1".
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$ftr.Range.Find.Execute("5", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2) | Out-Null
